# Refresh the IFRS financial figures for 테이팩스 (rows 2-9, columns D:AJ).
# Mirrors an upstream data re-pull: most numeric columns get new values
# and a handful of now-unused columns (J, O, Y, Z, AD, AH, ...) are cleared
# per row, matching the shrinking "spans" as later rows lose more columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1115
$ws.Range("E2").Value = 141
$ws.Range("F2").Value = 141
$ws.Range("G2").Value = 157
$ws.Range("H2").Value = 125
$ws.Range("I2").Value = 125
$ws.Range("K2").Value = 1194
$ws.Range("L2").Value = 261
$ws.Range("M2").Value = 932
$ws.Range("N2").Value = 932
$ws.Range("P2").Value = 36
$ws.Range("Q2").Value = 165
$ws.Range("R2").Value = -99
$ws.Range("S2").Value = -7
$ws.Range("T2").Value = 33
$ws.Range("U2").Value = 132
$ws.Range("V2").Value = 6
$ws.Range("W2").Value = 12.67
$ws.Range("X2").Value = 11.21
$ws.Range("AA2").Value = 28.02
$ws.Range("AB2").Value = 2522.98
$ws.Range("AC2").Value = 1761
$ws.Range("AE2").Value = 13131
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 7100000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3
$ws.Range("D3").Value = 1133
$ws.Range("E3").Value = 111
$ws.Range("F3").Value = 111
$ws.Range("G3").Value = 111
$ws.Range("H3").Value = 95
$ws.Range("I3").Value = 95
$ws.Range("K3").Value = 1114
$ws.Range("L3").Value = 536
$ws.Range("M3").Value = 578
$ws.Range("N3").Value = 578
$ws.Range("P3").Value = 36
$ws.Range("Q3").Value = 135
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = -151
$ws.Range("T3").Value = 48
$ws.Range("U3").Value = 87
$ws.Range("V3").Value = 304
$ws.Range("W3").Value = 9.83
$ws.Range("X3").Value = 8.41
$ws.Range("Y3").Value = 12.62
$ws.Range("Z3").Value = 8.26
$ws.Range("AA3").Value = 92.64
$ws.Range("AB3").Value = 1529.26
$ws.Range("AC3").Value = 1802
$ws.Range("AE3").Value = 13525
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 4274000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()

# Row 4
$ws.Range("D4").Value = 636
$ws.Range("E4").Value = 53
$ws.Range("F4").Value = 53
$ws.Range("G4").Value = 31
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("K4").Value = 1644
$ws.Range("L4").Value = 930
$ws.Range("M4").Value = 714
$ws.Range("N4").Value = 714
$ws.Range("P4").Value = 36
$ws.Range("Q4").Value = 54
$ws.Range("R4").Value = -1029
$ws.Range("S4").Value = 1008
$ws.Range("T4").Value = 13
$ws.Range("U4").Value = 40
$ws.Range("V4").Value = 613
$ws.Range("W4").Value = 8.4
$ws.Range("X4").Value = 0.77
$ws.Range("Y4").Value = 0.76
$ws.Range("Z4").Value = 0.36
$ws.Range("AA4").Value = 130.33
$ws.Range("AB4").Value = 1908.76
$ws.Range("AC4").Value = 115
$ws.Range("AE4").Value = 16704
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 4274000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("AD4").ClearContents()
$ws.Range("AH4").ClearContents()

# Row 5
$ws.Range("D5").Value = 1059
$ws.Range("E5").Value = 120
$ws.Range("F5").Value = 120
$ws.Range("G5").Value = 95
$ws.Range("H5").Value = 76
$ws.Range("I5").Value = 76
$ws.Range("K5").Value = 1598
$ws.Range("L5").Value = 718
$ws.Range("M5").Value = 880
$ws.Range("N5").Value = 880
$ws.Range("P5").Value = 38
$ws.Range("Q5").Value = 122
$ws.Range("R5").Value = -38
$ws.Range("S5").Value = -103
$ws.Range("T5").Value = 16
$ws.Range("U5").Value = 106
$ws.Range("V5").Value = 424
$ws.Range("W5").Value = 11.29
$ws.Range("X5").Value = 7.23
$ws.Range("Y5").Value = 9.6
$ws.Range("Z5").Value = 4.72
$ws.Range("AA5").Value = 81.51000000000001
$ws.Range("AB5").Value = 2233.21
$ws.Range("AC5").Value = 1747
$ws.Range("AD5").Value = 18.32
$ws.Range("AE5").Value = 18660
$ws.Range("AF5").Value = 1.71
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 4717650
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6
$ws.Range("D6").Value = 1085
$ws.Range("E6").Value = 75
$ws.Range("F6").Value = 75
$ws.Range("G6").Value = 67
$ws.Range("H6").Value = 51
$ws.Range("I6").Value = 51
$ws.Range("K6").Value = 1602
$ws.Range("L6").Value = 671
$ws.Range("M6").Value = 931
$ws.Range("N6").Value = 931
$ws.Range("P6").Value = 38
$ws.Range("Q6").Value = 45
$ws.Range("R6").Value = -9
$ws.Range("S6").Value = -36
$ws.Range("T6").Value = 12
$ws.Range("U6").Value = 33
$ws.Range("V6").Value = 385
$ws.Range("W6").Value = 6.94
$ws.Range("X6").Value = 4.7
$ws.Range("Y6").Value = 5.63
$ws.Range("Z6").Value = 3.19
$ws.Range("AA6").Value = 72.05
$ws.Range("AB6").Value = 2367.84
$ws.Range("AC6").Value = 1081
$ws.Range("AD6").Value = 19.47
$ws.Range("AE6").Value = 19740
$ws.Range("AF6").Value = 1.07
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 4717650
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 1140
$ws.Range("E7").Value = 77
$ws.Range("I7").Value = 57
$ws.Range("W7").Value = 6.75
$ws.Range("AC7").Value = 1208
$ws.Range("AD7").Value = 16.72
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 1295
$ws.Range("E8").Value = 102
$ws.Range("I8").Value = 89
$ws.Range("W8").Value = 7.88
$ws.Range("AC8").Value = 1887
$ws.Range("AD8").Value = 10.71
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
